$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ31510755",
    "summ31589804",
    "summ31672444",
    "summ31742115",
    "summ31821929",
    "summ31889704",
    "summ31973349",
    "summ32056012",
    "summ32122929",
    "summ32212964",
    "summ32290089",
    "summ32356641",
    "summ32444108",
    "summ32522301",
    "summ32595536",
    "summ32675571"
)

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
}
